$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the whole "*Link to Billing Information" bullet paragraph.
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "Link to Billing Information") {
        $p.Range.Delete()
    }
}

# ---------------------------------------------------------------------------
# 2. Split the "Delete function" run into "Delete " / "profile " / "function"
#    and move the (hidden) _GoBack bookmark so it now sits between
#    "profile " and "function" (Word keeps only one bookmark per name, so
#    re-adding it here also removes it from its old location).
# ---------------------------------------------------------------------------
$findRange = $d.Content
$found = $findRange.Find.Execute("Delete function")

$paraStart = $findRange.Start          # start of "Delete function"
$origSize = $findRange.Font.Size       # preserve the run's current font size
$insertPos = $paraStart + 7            # right after "Delete "

# Insert the new word; this currently just extends the existing run.
$ip = $d.Range($insertPos, $insertPos)
$ip.InsertAfter("profile ")

# Force Word to split "Delete " / "profile " into their own runs by
# nudging then restoring the character formatting of the inserted text.
$splitRange = $d.Range($insertPos, $insertPos + 8)
$splitRange.Font.Size = $origSize + 1
$splitRange.Font.Size = $origSize

# Re-seat the _GoBack bookmark right after "profile " (i.e. immediately
# before "function"), which also removes it from its previous location.
$bmPos = $insertPos + 8
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
